$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "262.11"
Set-TextValue "E2" "0.49%"
Set-TextValue "D3" "26.66"
Set-TextValue "E3" "-1.99%"
Set-TextValue "D4" "4.699"
Set-TextValue "E4" "0.14%"
Set-TextValue "D5" "0.06077"
Set-TextValue "E5" "-0.81%"
Set-TextValue "D6" "6.703"
Set-TextValue "E6" "0.70%"
Set-TextValue "D7" "0.8509"
Set-TextValue "E7" "-0.27%"
Set-TextValue "D8" "0.9083"
Set-TextValue "E8" "-1.40%"
Set-TextValue "D9" "0.1401"
Set-TextValue "E9" "-0.52%"
Set-TextValue "D10" "0.05101"
Set-TextValue "E10" "9.52%"
Set-TextValue "D11" "0.07093"
Set-TextValue "E11" "0.05%"
Set-TextValue "D12" "0.03117"
Set-TextValue "E12" "1.99%"
Set-TextValue "D13" "0.09044"
Set-TextValue "E13" "-0.20%"
Set-TextValue "D14" "0.001533"
Set-TextValue "E14" "-0.56%"
Set-TextValue "D15" "0.0006174"
Set-TextValue "E15" "1.24%"
Set-TextValue "D16" "0.005982"
Set-TextValue "E16" "-0.99%"
Set-TextValue "E17" "-0.06%"
Set-TextValue "D18" "3.170"
Set-TextValue "E18" "0.80%"
Set-TextValue "D19" "2.146"
Set-TextValue "E19" "-0.79%"
Set-TextValue "E21" "-2.25%"
Set-TextValue "D22" "4.117"
Set-TextValue "E22" "0.83%"
Set-TextValue "D23" "0.04236"
Set-TextValue "E23" "-0.17%"
Set-TextValue "E24" "-3.03%"
Set-TextValue "D25" "0.004059"
Set-TextValue "E25" "6.81%"
Set-TextValue "E26" "0.05%"
Set-TextValue "E27" "23.04%"
Set-TextValue "D40" "0.03954"
Set-TextValue "E40" "2.11%"
Set-TextValue "D41" "0.1112"
Set-TextValue "E41" "-0.01%"
Set-TextValue "D42" "0.004189"
Set-TextValue "E42" "2.76%"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002061"
Set-TextValue "E43" "-7.03%"
$ws.Range("B44").Value = "LocalTraders"
$ws.Range("C44").Value = "https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct"
Set-TextValue "D44" "0.01295"
Set-TextValue "E44" "-20.77%"
Set-TextValue "D45" "0.00005121"
Set-TextValue "E45" "-0.67%"
Set-TextValue "E46" "0.05%"
Set-TextValue "D48" "0.2584"
Set-TextValue "E48" "90.72%"
Set-TextValue "E49" "0.05%"
Set-TextValue "E50" "0.05%"

Write-Output "Applied all cell updates"
